$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the image filenames referenced in the "Time Series" column (D)
# from the 2025-04-07 plot exports to the newer 2025-05-05 exports.
$ws.Range("D2").Value = "Commercial_LONGFINSQUID_Landings_LBS_2025-05-05.png"
$ws.Range("D3").Value = "N_Commercial_Vessels_Landing_LONGFINSQUID_2025-05-05.png"
$ws.Range("D4").Value = "TOTALANNUALREV_LONGFINSQUID_2023Dols_2025-05-05.png"

# Reflect the author's final selection/active cell position.
$ws.Range("D4").Select()
